# TEMPLATE-IMPOR-SISWA.xlsx update:
#  - add a new "Petunjuk" (instructions) sheet right after Sheet1
#  - protect the Sheet1 data area, leaving A2:Q12000 ("data siswa") editable
#  - update the remembered selections on both sheets

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- new "Petunjuk" sheet, placed right after Sheet1 --------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Petunjuk"
$ws2.Tab.Color = 65535   # OLE (BGR) 00FFFF00 -> RGB FFFF00 (yellow)

$ws2.Range("A1").Value = "Jenis Kelamin cukup isi dengan L atau P"
$ws2.Range("A2").Value = "Tanggal lahir harus dalam format YYYY-MM-DD"
$ws2.Range("A3").Value = "No induk diisi dengan 9 karakter"
$ws2.Range("A4").Value = "NISN diisi dengan 10 karakter"

$ws2.Columns.Item(1).ColumnWidth = 54.33

# keep a (harmless, content-free) record of the sheet's former outline
# depth, matching the template's sheetFormatPr, without leaving a stray row
$ws2.Rows.Item(5).OutlineLevel = 3
$ws2.Rows.Item(5).Delete()

# force an explicit (empty) headerFooter section, like the rest of the workbook
$ws2.PageSetup.CenterHeader = ""

# --- protect the student-data sheet, keep the data range editable -------
$ws1.Protection.AllowEditRanges.Add("data siswa", $ws1.Range("A2:Q12000"))
$ws1.Protect("data siswa")

# --- restore the remembered selections -----------------------------------
$ws1.Range("J11").Select()
$ws2.Range("A13").Select()
